$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: Fecha = "08/01/25" (text, like the existing date-like entries
# in column A), Valor = 1 (number).
#
# A plain `.Value = "08/01/25"` gets auto-sniffed by Excel as a date literal
# and stored as a date serial number, which doesn't match the source data
# (a literal text string). Forcing the cell to the "@" (Text) number format
# before assigning the value keeps it as literal text; ClearFormats()
# afterwards drops the explicit "Text" style again so the cell ends up with
# no style override, exactly like the other data rows (A2:B4).
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "08/01/25"
$ws.Range("A5").ClearFormats()

$ws.Range("B5").Value = 1
